# fix: added emp id in sheet to manage the data
#
# Inserts a new "EmpId" column in front of the existing Employee/Paid
# columns, numbers the four existing rows 1-4, and bumps Alice's paid
# amount from 3900 to 3400.56. Also reproduces the formatting tweaks
# that came along with the edit (bold header/data cells, column width
# for the Employee column, row heights, and the new selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A; this shifts the existing "Employee" column to B
# and "Paid" to C (carrying their values/styles with them).
$ws.Columns("A").Insert()

# New header + sequential employee ids in column A.
$ws.Range("A1").Value = "EmpId"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

# Update Alice's paid amount.
$ws.Range("C2").Value = 3400.56

# A1 becomes bold, but (unlike the B1/C1 header cells) keeps no border
# and plain general/bottom alignment.
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Borders.LineStyle = -4142
$ws.Range("A1").HorizontalAlignment = 1
$ws.Range("A1").VerticalAlignment = -4107

# Widen the Employee column and normalise row heights.
$ws.Columns("B").ColumnWidth = 14.11
$ws.Rows("2:5").RowHeight = 13.8

# Matches the author's final selection.
$ws.Range("D11").Select()
